$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Content.Paragraphs.Item($index)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# --- 1. Merge split runs (no visible text change) ---
Set-ParaText 22 "Search Bar in the middle where location is inputted to help direct google maps in direction of user."
Set-ParaText 30 "Able to See Dive Sites on the map"
Set-ParaText 34 "Has a more details page which redirects you to a more detailed information page on the site"
Set-ParaText 50 "Able to click anywhere and add new dive sites."

# --- 2. DiveSiteDB table field-name edits (indices measured before any table mutation) ---
Set-ParaText 102 "Latitude (*)"
Set-ParaText 105 "longitude(*)"
Set-ParaText 108 "siteName(*)"
Set-ParaText 111 "areaName(*)"
Set-ParaText 114 "Description(*)"
Set-ParaText 133 "API to give current weather/future weather"
Set-ParaText 135 "siteType(*)"

# --- 3. Add "Bio" row to the User DB table (do structural edit last) ---
$userTable = $d.Tables.Item(1)
$newRow = $userTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Bio"
$newRow.Cells.Item(2).Range.Text = "Tech Freak, loves the outdoors"

Write-Output "done"
